$wb = $excel.ActiveWorkbook

# --- Sheet "Variables" (sheet1): add two new variable rows (41, 42) ---
$wsVar = $wb.Worksheets.Item("Variables")

# Row 41: inccanc / incident first occuring cancer / integer
$wsVar.Range("B41").Value = "inccanc"
$wsVar.Range("C41").Value = "incident first occuring cancer"
$wsVar.Range("D41").Value = "integer"

# Row 42: dcens_canc / censored age for cancer at FUP5 / date
$wsVar.Range("B42").Value = "dcens_canc"
$wsVar.Range("C42").Value = "censored age for cancer at FUP5"
$wsVar.Range("D42").Value = "date"

# Highlight the two new rows with a yellow fill (matches the rest of the row styling)
$wsVar.Range("B41:D42").Interior.Color = 65535

# --- Sheet "Categories" (sheet2): add category rows for the new "inccanc" variable ---
$wsCat = $wb.Worksheets.Item("Categories")

# Row 73: inccanc = 0 -> No
$wsCat.Range("A73").Value = "inccanc"
$wsCat.Range("B73").Value = 0
$wsCat.Range("C73").Value = "No"

# Row 74: inccanc = 1 -> Yes
$wsCat.Range("A74").Value = "inccanc"
$wsCat.Range("B74").Value = 1
$wsCat.Range("C74").Value = "Yes"

# Column A on these rows reuses the same (non-themed) font already used elsewhere
# in the sheet for variable-name cells (e.g. A49), then gets the yellow highlight.
$wsCat.Range("A49").Copy() | Out-Null
$wsCat.Range("A73:A74").PasteSpecial(-4122) | Out-Null
$wsCat.Range("A73:C74").Interior.Color = 65535

$excel.CutCopyMode = 0

# --- Restore selections to match where the author ended up ---
$wsVar.Activate()
$wsVar.Range("F35").Select()

$wsCat.Activate()
$wsCat.Range("C57").Select()

$wsVar.Activate()
